$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values look numeric,
# so Excel keeps them as text (matching the original inline-string cells).
$textFormatCells = @("D5", "D6", "D10", "D11", "D15", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values row by row
$ws.Range("D2").Value = '57.295.54'
$ws.Range("E2").Value = '  +1.94%  '

$ws.Range("D3").Value = '3.265.67'
$ws.Range("E3").Value = '  +1.15%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '398.29'
$ws.Range("E5").Value = '  +0.30%  '

$ws.Range("D6").Value = '108.86'
$ws.Range("E6").Value = '  -1.85%  '

$ws.Range("E7").Value = '  +4.60%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  +0.19%  '

$ws.Range("D10").Value = '39.21'
$ws.Range("E10").Value = '  -0.24%  '

$ws.Range("D11").Value = '0.0953'
$ws.Range("E11").Value = '  +4.41%  '

$ws.Range("E12").Value = '  +1.14%  '

$ws.Range("D13").Value = '3.779.26'
$ws.Range("E13").Value = '  +1.14%  '

$ws.Range("E14").Value = '  +2.21%  '

$ws.Range("D15").Value = '18.95'
$ws.Range("E15").Value = '  -0.50%  '

$ws.Range("D16").Value = '3.263.34'
$ws.Range("E16").Value = '  +1.04%  '

$ws.Range("E17").Value = '  -1.54%  '

$ws.Range("E18").Value = '  +3.25%  '

$ws.Range("D19").Value = '57.111.22'
$ws.Range("E19").Value = '  +1.91%  '

$ws.Range("D20").Value = '3.31'
$ws.Range("E20").Value = '  +0.01%  '

$ws.Range("E21").Value = '  +4.54%  '

$ws.Range("D22").Value = '12.94'
$ws.Range("E22").Value = '  -0.43%  '

$ws.Range("D23").Value = '297.14'
$ws.Range("E23").Value = '  -0.66%  '

$ws.Range("D24").Value = '74.02'
$ws.Range("E24").Value = '  -1.83%  '

$ws.Range("D25").Value = '3.19'
$ws.Range("E25").Value = '  -1.28%  '

$ws.Range("D26").Value = '28.13'

$ws.Range("D27").Value = '4.39'
$ws.Range("E27").Value = '  +0.17%  '

$ws.Range("D28").Value = '7.86'
$ws.Range("E28").Value = '  -4.19%  '

$ws.Range("E29").Value = '  -0.82%  '

$ws.Range("D30").Value = '0.169'
$ws.Range("E30").Value = '  -2.07%  '

$ws.Range("E31").Value = '  -0.40%  '

$ws.Range("E32").Value = '  +1.83%  '

$ws.Range("D33").Value = '11.19'
$ws.Range("E33").Value = '  +0.52%  '

$ws.Range("D34").Value = '40.05'
$ws.Range("E34").Value = '  +10.68%  '

$ws.Range("D35").Value = '0.0496'
$ws.Range("E35").Value = '  +0.70%  '

$ws.Range("E36").Value = '  +1.10%  '

$ws.Range("D37").Value = '51.40'
$ws.Range("E37").Value = '  +0.21%  '

$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.02%  '

$ws.Range("E39").Value = '  -1.78%  '

$ws.Range("D40").Value = '3.05'
$ws.Range("E40").Value = '  -3.25%  '

$ws.Range("D41").Value = '137.36'
$ws.Range("E41").Value = '  +1.71%  '

$ws.Range("E42").Value = '  +1.72%  '

$ws.Range("D43").Value = '0.287'
$ws.Range("E43").Value = '  +1.61%  '

$ws.Range("D44").Value = '1.88'
$ws.Range("E44").Value = '  -2.29%  '

$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = '16.78'
$ws.Range("E45").Value = '  -3.05%  '

$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '3.88'
$ws.Range("E46").Value = '  -3.46%  '

$ws.Range("D47").Value = '22.48'
$ws.Range("E47").Value = '  +0.89%  '

$ws.Range("E48").Value = '  +4.24%  '

$ws.Range("D49").Value = '2.151.45'
$ws.Range("E49").Value = '  +0.89%  '

$ws.Range("E50").Value = '  -0.12%  '

$ws.Range("D51").Value = '1.96'
$ws.Range("E51").Value = '  -7.99%  '
